# Auto-generated edit script applying the Gilgamesh_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across 8 sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 2749.8
$ws.Range("I28").Value = 949
$ws.Range("K28").Value = 949
$ws.Range("M28").Value = -464
$ws.Range("H100").Value = 4499.25
$ws.Range("J100").Value = 9999
$ws.Range("L100").Value = 9999
$ws.Range("N100").Value = -11081
$ws.Range("H116").Value = 9290.125
$ws.Range("I116").Value = 2437
$ws.Range("J116").Value = 9817.288
$ws.Range("K116").Value = 2437
$ws.Range("L116").Value = 9817.288
$ws.Range("M116").Value = 1005
$ws.Range("N116").Value = -16701.288
$ws.Range("H138").Value = 350278.5
$ws.Range("I138").Value = 3159.8215
$ws.Range("J138").Value = 571172.25
$ws.Range("K138").Value = 9479.4645
$ws.Range("L138").Value = 1713516.75
$ws.Range("M138").Value = -4339.4645
$ws.Range("N138").Value = -1723796.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3759.9363
$ws.Range("I32").Value = 3733
$ws.Range("J32").Value = 4999
$ws.Range("K32").Value = 3733
$ws.Range("L32").Value = 4999
$ws.Range("M32").Value = -3446
$ws.Range("N32").Value = -5573
$ws.Range("H45").Value = 50991.777
$ws.Range("I45").Value = 86881.8
$ws.Range("J45").Value = 6129.25
$ws.Range("K45").Value = 86881.8
$ws.Range("L45").Value = 6129.25
$ws.Range("M45").Value = -86504.8
$ws.Range("N45").Value = -6883.25
$ws.Range("H61").Value = 4081.303
$ws.Range("I61").Value = 2699.7917
$ws.Range("K61").Value = 2699.7917
$ws.Range("M61").Value = -2487.7917
$ws.Range("H74").Value = 267377.47
$ws.Range("I74").Value = 557261.4
$ws.Range("K74").Value = 557261.4
$ws.Range("M74").Value = -556387.4
$ws.Range("H77").Value = 267377.47
$ws.Range("I77").Value = 557261.4
$ws.Range("K77").Value = 2786307
$ws.Range("M77").Value = -2781939
$ws.Range("H110").Value = 2783
$ws.Range("I110").Value = 1472.3636
$ws.Range("K110").Value = 1472.3636
$ws.Range("M110").Value = 572.6364000000001
$ws.Range("H122").Value = 3992.923
$ws.Range("I122").Value = 3852.3635
$ws.Range("J122").Value = 4766
$ws.Range("K122").Value = 11557.0905
$ws.Range("L122").Value = 14298
$ws.Range("M122").Value = -9107.0905
$ws.Range("N122").Value = -19198
$ws.Range("H136").Value = 4081.303
$ws.Range("I136").Value = 2699.7917
$ws.Range("K136").Value = 8099.375100000001
$ws.Range("M136").Value = -5549.375100000001
$ws.Range("H138").Value = 96661.336
$ws.Range("J138").Value = 96661.336
$ws.Range("L138").Value = 96661.336
$ws.Range("N138").Value = -106941.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2783.0588
$ws.Range("I86").Value = 2651.923
$ws.Range("K86").Value = 2651.923
$ws.Range("M86").Value = -1528.923
$ws.Range("H89").Value = 2783.0588
$ws.Range("I89").Value = 2651.923
$ws.Range("K89").Value = 13259.615
$ws.Range("M89").Value = -7643.614999999998
$ws.Range("H99").Value = 146794.42
$ws.Range("I99").Value = 202512.4
$ws.Range("K99").Value = 202512.4
$ws.Range("M99").Value = -201014.4
$ws.Range("H105").Value = 13002145
$ws.Range("I105").Value = 590128.5600000001
$ws.Range("K105").Value = 590128.5600000001
$ws.Range("M105").Value = -588381.5600000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4360.0435
$ws.Range("J31").Value = 5141.769
$ws.Range("L31").Value = 5141.769
$ws.Range("N31").Value = -5731.769
$ws.Range("H34").Value = 4360.0435
$ws.Range("J34").Value = 5141.769
$ws.Range("L34").Value = 5141.769
$ws.Range("N34").Value = -5545.769
$ws.Range("H99").Value = 9151.333000000001
$ws.Range("I99").Value = 15315.833
$ws.Range("J99").Value = 5041.6665
$ws.Range("K99").Value = 15315.833
$ws.Range("L99").Value = 5041.6665
$ws.Range("M99").Value = -13817.833
$ws.Range("N99").Value = -8037.6665
$ws.Range("H122").Value = 4205.7144
$ws.Range("J122").Value = 5133.2856
$ws.Range("L122").Value = 15399.8568
$ws.Range("N122").Value = -20299.8568
$ws.Range("H126").Value = 9151.333000000001
$ws.Range("I126").Value = 15315.833
$ws.Range("J126").Value = 5041.6665
$ws.Range("K126").Value = 45947.499
$ws.Range("L126").Value = 15124.9995
$ws.Range("M126").Value = -43477.499
$ws.Range("N126").Value = -20064.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 3910.7693
$ws.Range("J81").Value = 4844
$ws.Range("L81").Value = 14532
$ws.Range("N81").Value = -16778
$ws.Range("H84").Value = 3910.7693
$ws.Range("J84").Value = 4844
$ws.Range("L84").Value = 43596
$ws.Range("N84").Value = -54828
$ws.Range("H134").Value = 1392.9231
$ws.Range("I134").Value = 1392.9231
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4178.7693
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 891.2307000000001
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 12000
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 12000
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 36000
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -40900
$ws.Range("H132").Value = 2092.6428
$ws.Range("I132").Value = 1524.75
$ws.Range("J132").Value = 2319.8
$ws.Range("K132").Value = 4574.25
$ws.Range("L132").Value = 6959.400000000001
$ws.Range("M132").Value = -2044.25
$ws.Range("N132").Value = -12019.4
$ws.Range("H139").Value = 84461.375
$ws.Range("J139").Value = 86485
$ws.Range("L139").Value = 86485
$ws.Range("N139").Value = -96765

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2815.4546
$ws.Range("I7").Value = 2815.4546
$ws.Range("K7").Value = 2815.4546
$ws.Range("M7").Value = -2703.4546
$ws.Range("H40").Value = 34021.97
$ws.Range("I40").Value = 40620.82
$ws.Range("J40").Value = 3227.3333
$ws.Range("K40").Value = 40620.82
$ws.Range("L40").Value = 3227.3333
$ws.Range("M40").Value = -40484.82
$ws.Range("N40").Value = -3499.3333
$ws.Range("H126").Value = 2815.4546
$ws.Range("I126").Value = 2815.4546
$ws.Range("K126").Value = 8446.363799999999
$ws.Range("M126").Value = -5976.363799999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 19225.834
$ws.Range("J41").Value = 19296.4
$ws.Range("L41").Value = 19296.4
$ws.Range("N41").Value = -20076.4
$ws.Range("H107").Value = 715.7273
$ws.Range("I107").Value = 403.14285
$ws.Range("K107").Value = 1209.42855
$ws.Range("M107").Value = 710.5714499999999
$ws.Range("H122").Value = 12502189
$ws.Range("I122").Value = 2285.4167
$ws.Range("J122").Value = 31252044
$ws.Range("K122").Value = 6856.250100000001
$ws.Range("L122").Value = 93756132
$ws.Range("M122").Value = -4406.250100000001
$ws.Range("N122").Value = -93761032
$ws.Range("H126").Value = 1149
$ws.Range("I126").Value = 1149
$ws.Range("K126").Value = 3447
$ws.Range("M126").Value = -977
$ws.Range("H132").Value = 4701.933
$ws.Range("I132").Value = 5504
$ws.Range("K132").Value = 16512
$ws.Range("M132").Value = -13982
